# The presentation currently ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" colours (wired to the Notes Master)
#   ppt/theme/theme2.xml -> "Integral" colours      (wired to the Slide Master / Design)
#
# The authored revision swaps the colour palettes held by those two theme
# parts: the theme that drives the Slide Master ends up using the plain
# "Office" palette, while the other part ends up holding the "Integral"
# palette it used to replace. The font scheme and format scheme are
# identical between the two parts already, so the only real difference is
# the 12-colour scheme (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# PowerPoint's object model exposes exactly this surface for editing a
# design's theme: ActivePresentation.SlideMaster.Theme.ThemeColorScheme,
# a 12-item indexed collection of ColorFormat objects (RGB settable). Drive
# the active design's palette from "Integral" to "Office" through it.

function ToComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$colors = $p.SlideMaster.Theme.ThemeColorScheme

# MsoThemeColorSchemeIndex order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 0; $i -lt $officeColors.Count; $i++) {
    $colors.Item($i + 1).RGB = ToComRGB $officeColors[$i]
}
